$wb = $excel.ActiveWorkbook

# --- Sheet 1: semantic_aspect_model_schema ---
$ws1 = $wb.Worksheets.Item("semantic_aspect_model_schema")

# Column width adjustments (widen columns to account for longer double-underscore field names)
$ws1.Columns.Item(7).ColumnWidth = 24.333333333333332
$ws1.Columns.Item(15).ColumnWidth = 28.0
$ws1.Columns.Item(17).ColumnWidth = 28.0
$ws1.Columns.Item(20).ColumnWidth = 19.5
$ws1.Columns.Item(21).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(22).ColumnWidth = 36.333333333333336
$ws1.Columns.Item(23).ColumnWidth = 34.0
$ws1.Columns.Item(24).ColumnWidth = 38.833333333333336
$ws1.Columns.Item(25).ColumnWidth = 38.833333333333336
$ws1.Columns.Item(26).ColumnWidth = 37.5
$ws1.Columns.Item(27).ColumnWidth = 24.333333333333332
$ws1.Columns.Item(28).ColumnWidth = 36.333333333333336
$ws1.Columns.Item(29).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(30).ColumnWidth = 26.833333333333332
$ws1.Columns.Item(35).ColumnWidth = 43.5
$ws1.Columns.Item(36).ColumnWidth = 47.166666666666664
$ws1.Columns.Item(37).ColumnWidth = 46.0
$ws1.Columns.Item(38).ColumnWidth = 24.333333333333332
$ws1.Columns.Item(40).ColumnWidth = 46.0
$ws1.Columns.Item(41).ColumnWidth = 47.166666666666664
$ws1.Columns.Item(42).ColumnWidth = 41.166666666666664
$ws1.Columns.Item(43).ColumnWidth = 46.0
$ws1.Columns.Item(44).ColumnWidth = 46.0
$ws1.Columns.Item(45).ColumnWidth = 44.833333333333336
$ws1.Columns.Item(46).ColumnWidth = 47.166666666666664
$ws1.Columns.Item(47).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(48).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(49).ColumnWidth = 26.833333333333332
$ws1.Columns.Item(50).ColumnWidth = 47.166666666666664
$ws1.Columns.Item(51).ColumnWidth = 34.0
$ws1.Columns.Item(52).ColumnWidth = 24.333333333333332
$ws1.Columns.Item(53).ColumnWidth = 30.333333333333332
$ws1.Columns.Item(54).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(55).ColumnWidth = 41.166666666666664
$ws1.Columns.Item(56).ColumnWidth = 49.166666666666664
$ws1.Columns.Item(57).ColumnWidth = 49.166666666666664
$ws1.Columns.Item(58).ColumnWidth = 47.166666666666664
$ws1.Columns.Item(63).ColumnWidth = 26.833333333333332
$ws1.Columns.Item(64).ColumnWidth = 38.833333333333336
$ws1.Columns.Item(65).ColumnWidth = 30.333333333333332

# Header text updates (row 1): single underscore -> double underscore field separator
$ws1.Range("G1").Value = "precedingPfIds[0]__id"
$ws1.Range("O1").Value = "companyIds[0]__companyId"
$ws1.Range("Q1").Value = "productIds[0]__productId"
$ws1.Range("T1").Value = "pcf__declaredUnit"
$ws1.Range("U1").Value = "pcf__unitaryProductAmount"
$ws1.Range("V1").Value = "pcf__productMassPerDeclaredUnit"
$ws1.Range("W1").Value = "pcf__exemptedEmissionsPercent"
$ws1.Range("X1").Value = "pcf__exemptedEmissionsDescription"
$ws1.Range("Y1").Value = "pcf__boundaryProcessesDescription"
$ws1.Range("Z1").Value = "pcf__geographyCountrySubdivision"
$ws1.Range("AA1").Value = "pcf__geographyCountry"
$ws1.Range("AB1").Value = "pcf__geographyRegionOrSubregion"
$ws1.Range("AC1").Value = "pcf__referencePeriodStart"
$ws1.Range("AD1").Value = "pcf__referencePeriodEnd"
$ws1.Range("AE1").Value = "pcf__crossSectoralStandardsUsed[0]__crossSectoralStandard"
$ws1.Range("AF1").Value = "pcf__productOrSectorSpecificRules[0]__extWBCSD_operator"
$ws1.Range("AG1").Value = "pcf__productOrSectorSpecificRules[0]__productOrSectorSpecificRules[0]__ruleName"
$ws1.Range("AH1").Value = "pcf__productOrSectorSpecificRules[0]__extWBCSD_otherOperatorName"
$ws1.Range("AI1").Value = "pcf__extWBCSD_characterizationFactors"
$ws1.Range("AJ1").Value = "pcf__extWBCSD_allocationRulesDescription"
$ws1.Range("AK1").Value = "pcf__extTFS_allocationWasteIncineration"
$ws1.Range("AL1").Value = "pcf__primaryDataShare"
$ws1.Range("AM1").Value = "pcf__secondaryEmissionFactorSources[0]__secondaryEmissionFactorSource"
$ws1.Range("AN1").Value = "pcf__dataQualityRating__coveragePercent"
$ws1.Range("AO1").Value = "pcf__dataQualityRating__technologicalDQR"
$ws1.Range("AP1").Value = "pcf__dataQualityRating__temporalDQR"
$ws1.Range("AQ1").Value = "pcf__dataQualityRating__geographicalDQR"
$ws1.Range("AR1").Value = "pcf__dataQualityRating__completenessDQR"
$ws1.Range("AS1").Value = "pcf__dataQualityRating__reliabilityDQR"
$ws1.Range("AT1").Value = "pcf__extWBCSD_packagingEmissionsIncluded"
$ws1.Range("AU1").Value = "pcf__pcfExcludingBiogenic"
$ws1.Range("AV1").Value = "pcf__pcfIncludingBiogenic"
$ws1.Range("AW1").Value = "pcf__fossilGhgEmissions"
$ws1.Range("AX1").Value = "pcf__biogenicCarbonEmissionsOtherThanCO2"
$ws1.Range("AY1").Value = "pcf__biogenicCarbonWithdrawal"
$ws1.Range("AZ1").Value = "pcf__dlucGhgEmissions"
$ws1.Range("BA1").Value = "pcf__extTFS_luGhgEmissions"
$ws1.Range("BB1").Value = "pcf__aircraftGhgEmissions"
$ws1.Range("BC1").Value = "pcf__extWBCSD_packagingGhgEmissions"
$ws1.Range("BD1").Value = "pcf__distributionStagePcfExcludingBiogenic"
$ws1.Range("BE1").Value = "pcf__distributionStagePcfIncludingBiogenic"
$ws1.Range("BF1").Value = "pcf__distributionStageFossilGhgEmissions"
$ws1.Range("BG1").Value = "pcf__distributionStageBiogenicCarbonEmissionsOtherThanCO2"
$ws1.Range("BH1").Value = "pcf__distributionStageBiogenicCarbonWithdrawal"
$ws1.Range("BI1").Value = "pcf__extTFS_distributionStageDlucGhgEmissions"
$ws1.Range("BJ1").Value = "pcf__extTFS_distributionStageLuGhgEmissions"
$ws1.Range("BK1").Value = "pcf__carbonContentTotal"
$ws1.Range("BL1").Value = "pcf__extWBCSD_fossilCarbonContent"
$ws1.Range("BM1").Value = "pcf__carbonContentBiogenic"

# --- Sheet 2: description ---
$ws2 = $wb.Worksheets.Item("description")

$ws2.Range("A3").Value = "1. Columns highlighted in olive green are digital twin fields."
$ws2.Range("B5").Value = "Digital Twin Field Name: id"
$ws2.Range("B6").Value = "Digital Twin Field Name: manufacturerPartId"
$ws2.Range("B7").Value = "Digital Twin Field Name: digitalTwinType"
$ws2.Range("A11").Value = "precedingPfIds[0]__id"
$ws2.Range("A19").Value = "companyIds[0]__companyId"
$ws2.Range("A21").Value = "productIds[0]__productId"
$ws2.Range("A24").Value = "pcf__declaredUnit"
$ws2.Range("A25").Value = "pcf__unitaryProductAmount"
$ws2.Range("A26").Value = "pcf__productMassPerDeclaredUnit"
$ws2.Range("A27").Value = "pcf__exemptedEmissionsPercent"
$ws2.Range("A28").Value = "pcf__exemptedEmissionsDescription"
$ws2.Range("A29").Value = "pcf__boundaryProcessesDescription"
$ws2.Range("A30").Value = "pcf__geographyCountrySubdivision"
$ws2.Range("A31").Value = "pcf__geographyCountry"
$ws2.Range("A32").Value = "pcf__geographyRegionOrSubregion"
$ws2.Range("A33").Value = "pcf__referencePeriodStart"
$ws2.Range("A34").Value = "pcf__referencePeriodEnd"
$ws2.Range("A35").Value = "pcf__crossSectoralStandardsUsed[0]__crossSectoralStandard"
$ws2.Range("A36").Value = "pcf__productOrSectorSpecificRules[0]__extWBCSD_operator"
$ws2.Range("A37").Value = "pcf__productOrSectorSpecificRules[0]__productOrSectorSpecificRules[0]__ruleName"
$ws2.Range("A38").Value = "pcf__productOrSectorSpecificRules[0]__extWBCSD_otherOperatorName"
$ws2.Range("A39").Value = "pcf__extWBCSD_characterizationFactors"
$ws2.Range("A40").Value = "pcf__extWBCSD_allocationRulesDescription"
$ws2.Range("A41").Value = "pcf__extTFS_allocationWasteIncineration"
$ws2.Range("A42").Value = "pcf__primaryDataShare"
$ws2.Range("A43").Value = "pcf__secondaryEmissionFactorSources[0]__secondaryEmissionFactorSource"
$ws2.Range("A44").Value = "pcf__dataQualityRating__coveragePercent"
$ws2.Range("A45").Value = "pcf__dataQualityRating__technologicalDQR"
$ws2.Range("A46").Value = "pcf__dataQualityRating__temporalDQR"
$ws2.Range("A47").Value = "pcf__dataQualityRating__geographicalDQR"
$ws2.Range("A48").Value = "pcf__dataQualityRating__completenessDQR"
$ws2.Range("A49").Value = "pcf__dataQualityRating__reliabilityDQR"
$ws2.Range("A50").Value = "pcf__extWBCSD_packagingEmissionsIncluded"
$ws2.Range("A51").Value = "pcf__pcfExcludingBiogenic"
$ws2.Range("A52").Value = "pcf__pcfIncludingBiogenic"
$ws2.Range("A53").Value = "pcf__fossilGhgEmissions"
$ws2.Range("A54").Value = "pcf__biogenicCarbonEmissionsOtherThanCO2"
$ws2.Range("A55").Value = "pcf__biogenicCarbonWithdrawal"
$ws2.Range("A56").Value = "pcf__dlucGhgEmissions"
$ws2.Range("A57").Value = "pcf__extTFS_luGhgEmissions"
$ws2.Range("A58").Value = "pcf__aircraftGhgEmissions"
$ws2.Range("A59").Value = "pcf__extWBCSD_packagingGhgEmissions"
$ws2.Range("A60").Value = "pcf__distributionStagePcfExcludingBiogenic"
$ws2.Range("A61").Value = "pcf__distributionStagePcfIncludingBiogenic"
$ws2.Range("A62").Value = "pcf__distributionStageFossilGhgEmissions"
$ws2.Range("A63").Value = "pcf__distributionStageBiogenicCarbonEmissionsOtherThanCO2"
$ws2.Range("A64").Value = "pcf__distributionStageBiogenicCarbonWithdrawal"
$ws2.Range("A65").Value = "pcf__extTFS_distributionStageDlucGhgEmissions"
$ws2.Range("A66").Value = "pcf__extTFS_distributionStageLuGhgEmissions"
$ws2.Range("A67").Value = "pcf__carbonContentTotal"
$ws2.Range("A68").Value = "pcf__extWBCSD_fossilCarbonContent"
$ws2.Range("A69").Value = "pcf__carbonContentBiogenic"
